$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; New = "63÷7=" },
    @{ Row = 1;  Col = 2; New = "79÷7=" },
    @{ Row = 1;  Col = 3; New = "39÷5=" },
    @{ Row = 1;  Col = 4; New = "51÷4=" },
    @{ Row = 1;  Col = 5; New = "47÷6=" },

    @{ Row = 5;  Col = 1; New = "93÷8=" },
    @{ Row = 5;  Col = 2; New = "48÷7=" },
    @{ Row = 5;  Col = 3; New = "27÷8=" },
    @{ Row = 5;  Col = 4; New = "18÷3=" },
    @{ Row = 5;  Col = 5; New = "47÷3=" },

    @{ Row = 9;  Col = 1; New = "16÷6=" },
    @{ Row = 9;  Col = 2; New = "65÷8=" },
    @{ Row = 9;  Col = 3; New = "31÷2=" },
    @{ Row = 9;  Col = 4; New = "69÷3=" },
    @{ Row = 9;  Col = 5; New = "49÷4=" },

    @{ Row = 13; Col = 1; New = "80÷3=" },
    @{ Row = 13; Col = 2; New = "14÷2=" },
    @{ Row = 13; Col = 3; New = "99÷4=" },
    @{ Row = 13; Col = 4; New = "26÷2=" },
    @{ Row = 13; Col = 5; New = "53÷3=" },

    @{ Row = 17; Col = 1; New = "26÷5=" },
    @{ Row = 17; Col = 2; New = "11÷8=" },
    @{ Row = 17; Col = 3; New = "68÷4=" },
    @{ Row = 17; Col = 4; New = "64÷8=" },
    @{ Row = 17; Col = 5; New = "21÷8=" }
)

foreach ($rep in $replacements) {
    $cell = $tbl.Cell($rep.Row, $rep.Col)
    $rng = $cell.Range
    # Exclude the trailing cell-mark character so only the visible text is replaced.
    $rng.End = $rng.End - 1
    $rng.Text = $rep.New
}
